$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 1 to hold the "Name"/"Marks" header
$ws.Rows.Item(1).Insert()

# Copy the existing header row's formatting (border + centered alignment)
# onto the new row, then re-tint just the fill color for the new header.
$ws.Range("A2:B2").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Fill in the new header text
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Marks"

# Give the new header row its own (lighter, accent-colored) fill
# (Accent 1 theme color, ~40% tint - lighter than the original header's fill)
$ws.Range("A1:B1").Interior.ThemeColor = 5
$ws.Range("A1:B1").Interior.TintAndShade = 0.39997558519241921


# Row heights: new header row + the row that used to be row 1
$ws.Rows.Item(1).RowHeight = 21
$ws.Rows.Item(2).RowHeight = 21

$ws.Range("D8").Select()
